$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '25.775.25'
$ws.Range('E2').Value = "'" + '  +0.34%  '
$ws.Range('D3').Value = "'" + '1.751.60'
$ws.Range('E3').Value = "'" + '  +0.47%  '
$ws.Range('D4').Value = "'" + '0.9993'
$ws.Range('E4').Value = "'" + '  -0.10%  '
$ws.Range('D5').Value = "'" + '236.90'
$ws.Range('E5').Value = "'" + '  -0.39%  '
$ws.Range('D6').Value = "'" + '0.9990'
$ws.Range('E6').Value = "'" + '  -0.12%  '
$ws.Range('D7').Value = "'" + '0.5063'
$ws.Range('E7').Value = "'" + '  +2.55%  '
$ws.Range('B8').Value = "'" + 'OKB'
$ws.Range('C8').Value = "'" + 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D8').Value = "'" + '40.84'
$ws.Range('E8').Value = "'" + '  -1.79%  '
$ws.Range('B9').Value = "'" + 'Cardano'
$ws.Range('C9').Value = "'" + 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').Value = "'" + '0.2713'
$ws.Range('E9').Value = "'" + '  +12.52%  '
$ws.Range('B10').Value = "'" + 'Dogecoin'
$ws.Range('C10').Value = "'" + 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Value = "'" + '0.06211'
$ws.Range('E10').Value = "'" + '  +4.05%  '
$ws.Range('B11').Value = "'" + 'WrappedEther'
$ws.Range('C11').Value = "'" + 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D11').Value = "'" + '1.751.74'
$ws.Range('E11').Value = "'" + '  +0.51%  '
$ws.Range('B12').Value = "'" + 'TRON'
$ws.Range('C12').Value = "'" + 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').Value = "'" + '0.06935'
$ws.Range('E12').Value = "'" + '  +1.35%  '
$ws.Range('B13').Value = "'" + 'Solana'
$ws.Range('C13').Value = "'" + 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D13').Value = "'" + '15.65'
$ws.Range('E13').Value = "'" + '  +6.55%  '
$ws.Range('B14').Value = "'" + 'Polygon'
$ws.Range('C14').Value = "'" + 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').Value = "'" + '0.6241'
$ws.Range('E14').Value = "'" + '  +7.27%  '
$ws.Range('B15').Value = "'" + 'Litecoin'
$ws.Range('C15').Value = "'" + 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D15').Value = "'" + '78.50'
$ws.Range('E15').Value = "'" + '  +1.74%  '
$ws.Range('B16').Value = "'" + 'Polkadot'
$ws.Range('C16').Value = "'" + 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D16').Value = "'" + '4.492'
$ws.Range('E16').Value = "'" + '  +0.69%  '
$ws.Range('B17').Value = "'" + 'BinanceUSD'
$ws.Range('C17').Value = "'" + 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D17').Value = "'" + '0.9988'
$ws.Range('B18').Value = "'" + 'Dai'
$ws.Range('C18').Value = "'" + 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D18').Value = "'" + '0.9994'
$ws.Range('E18').Value = "'" + '  -0.07%  '
$ws.Range('B19').Value = "'" + 'WrappedBTC'
$ws.Range('C19').Value = "'" + 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').Value = "'" + '25.798.73'
$ws.Range('E19').Value = "'" + '  +0.29%  '
$ws.Range('B20').Value = "'" + 'Avalanche'
$ws.Range('C20').Value = "'" + 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').Value = "'" + '11.71'
$ws.Range('E20').Value = "'" + '  +2.08%  '
$ws.Range('B21').Value = "'" + 'ShibaInu'
$ws.Range('C21').Value = "'" + 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').Value = "'" + '0.000006741'
$ws.Range('E21').Value = "'" + '  +4.67%  '
$ws.Range('B22').Value = "'" + 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = "'" + 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').Value = "'" + '1.973.71'
$ws.Range('E22').Value = "'" + '  +0.61%  '
$ws.Range('B23').Value = "'" + 'Uniswap'
$ws.Range('C23').Value = "'" + 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').Value = "'" + '4.057'
$ws.Range('E23').Value = "'" + '  +2.56%  '
$ws.Range('B24').Value = "'" + 'Cosmos'
$ws.Range('C24').Value = "'" + 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').Value = "'" + '8.264'
$ws.Range('E24').Value = "'" + '  +5.89%  '
$ws.Range('B25').Value = "'" + 'Chainlink'
$ws.Range('C25').Value = "'" + 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D25').Value = "'" + '5.194'
$ws.Range('E25').Value = "'" + '  +3.91%  '
$ws.Range('B26').Value = "'" + 'Monero'
$ws.Range('C26').Value = "'" + 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = "'" + '137.03'
$ws.Range('E26').Value = "'" + '  +0.81%  '
$ws.Range('D27').Value = "'" + '15.17'
$ws.Range('E27').Value = "'" + '  +4.54%  '
$ws.Range('B28').Value = "'" + 'Toncoin'
$ws.Range('C28').Value = "'" + 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').Value = "'" + '1.453'
$ws.Range('E28').Value = "'" + '  -1.26%  '
$ws.Range('B29').Value = "'" + 'LidoDAOToken'
$ws.Range('C29').Value = "'" + 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').Value = "'" + '1.798'
$ws.Range('E29').Value = "'" + '  -1.96%  '
$ws.Range('B30').Value = "'" + 'BitcoinCash'
$ws.Range('C30').Value = "'" + 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').Value = "'" + '103.02'
$ws.Range('E30').Value = "'" + '  +2.29%  '
$ws.Range('B31').Value = "'" + 'Stellar'
$ws.Range('C31').Value = "'" + 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').Value = "'" + '0.08283'
$ws.Range('E31').Value = "'" + '  +2.04%  '
$ws.Range('B32').Value = "'" + 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = "'" + 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = "'" + '3.747'
$ws.Range('E32').Value = "'" + '  -0.97%  '
$ws.Range('B33').Value = "'" + 'Filecoin'
$ws.Range('C33').Value = "'" + 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = "'" + '3.431'
$ws.Range('E33').Value = "'" + '  +2.12%  '
$ws.Range('B34').Value = "'" + 'Hedera'
$ws.Range('C34').Value = "'" + 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = "'" + '0.04392'
$ws.Range('E34').Value = "'" + '  +0.33%  '
$ws.Range('B35').Value = "'" + 'Frax'
$ws.Range('C35').Value = "'" + 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D35').Value = "'" + '0.9986'
$ws.Range('E35').Value = "'" + '  -0.12%  '
$ws.Range('B36').Value = "'" + 'HuobiToken'
$ws.Range('C36').Value = "'" + 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value = "'" + '2.646'
$ws.Range('E36').Value = "'" + '  +0.18%  '
$ws.Range('B37').Value = "'" + 'ARBITRUM'
$ws.Range('C37').Value = "'" + 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').Value = "'" + '1.009'
$ws.Range('E37').Value = "'" + '  -1.22%  '
$ws.Range('B38').Value = "'" + 'ImmutableX'
$ws.Range('C38').Value = "'" + 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').Value = "'" + '0.6050'
$ws.Range('E38').Value = "'" + '  -0.13%  '
$ws.Range('B39').Value = "'" + 'MXToken'
$ws.Range('C39').Value = "'" + 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').Value = "'" + '2.687'
$ws.Range('E39').Value = "'" + '  -0.78%  '
$ws.Range('B40').Value = "'" + 'RenderToken'
$ws.Range('C40').Value = "'" + 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').Value = "'" + '1.967'
$ws.Range('E40').Value = "'" + '  -5.22%  '
$ws.Range('B41').Value = "'" + 'VeChain'
$ws.Range('C41').Value = "'" + 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').Value = "'" + '0.01559'
$ws.Range('E41').Value = "'" + '  +5.00%  '
$ws.Range('B42').Value = "'" + 'PaxDollar'
$ws.Range('C42').Value = "'" + 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').Value = "'" + '0.9989'
$ws.Range('E42').Value = "'" + '  -0.13%  '
$ws.Range('B43').Value = "'" + 'Quant'
$ws.Range('C43').Value = "'" + 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').Value = "'" + '102.09'
$ws.Range('E43').Value = "'" + '  -0.86%  '
$ws.Range('B44').Value = "'" + 'TrustWalletToken'
$ws.Range('C44').Value = "'" + 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').Value = "'" + '0.7535'
$ws.Range('E44').Value = "'" + '  -3.20%  '
$ws.Range('B45').Value = "'" + 'TheSandbox'
$ws.Range('C45').Value = "'" + 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').Value = "'" + '0.3841'
$ws.Range('E45').Value = "'" + '  +2.06%  '
$ws.Range('B46').Value = "'" + 'FraxShare'
$ws.Range('C46').Value = "'" + 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').Value = "'" + '4.858'
$ws.Range('E46').Value = "'" + '  -5.33%  '
$ws.Range('B47').Value = "'" + 'Cronos'
$ws.Range('C47').Value = "'" + 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').Value = "'" + '0.05496'
$ws.Range('E47').Value = "'" + '  +7.56%  '
$ws.Range('B48').Value = "'" + 'Algorand'
$ws.Range('C48').Value = "'" + 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').Value = "'" + '0.1087'
$ws.Range('E48').Value = "'" + '  +1.75%  '
$ws.Range('D49').Value = "'" + '5.980'
$ws.Range('E49').Value = "'" + '  -0.10%  '
$ws.Range('B50').Value = "'" + 'Elrond'
$ws.Range('C50').Value = "'" + 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').Value = "'" + '30.31'
$ws.Range('E50').Value = "'" + '  +0.88%  '
$ws.Range('B51').Value = "'" + 'Aave'
$ws.Range('C51').Value = "'" + 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = "'" + '52.49'
$ws.Range('E51').Value = "'" + '  -0.39%  '
